# feat: add 2022-Q1 data
#
# 1) Insert a new sheet "2022-Q1" (fund-level detail) right before "总计",
#    matching the layout/style of the other quarterly sheets.
# 2) Insert a new top data row into "总计" with the 2022-Q1 roll-up figures,
#    pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" worksheet, positioned immediately before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q4Sheet = $wb.Worksheets.Item("2021-Q4")

$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# NOTE: sheet references resolve by position, so $totalSheet now points at
# the freshly-inserted "2022-Q1" sheet (it took over index 5). Re-fetch the
# real "总计" sheet by name now that it has shifted to the next index.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header styling (bold + border, s=2) from an existing quarter sheet
$q4Sheet.Range("A1:H1").Copy()
$q1.Range("A1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Copy the A-column "index" styling (bold + border, s=2) onto A2:A6
$q4Sheet.Range("A2").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)

# B, D, E, F, G look numeric ("014279", "4.95", ...) but must stay plain
# text, matching the inlineStr typing on the other quarter sheets - a bare
# leading apostrophe keeps Excel from silently coercing them to numbers
# (and from dropping the leading zero on fund codes) without touching the
# cell's formatting. C (fund name) is never numeric-looking, so it is set
# as-is. H (rank) is a genuine number on every quarter sheet.
$q1Data = @(
    @(0, "014279", "汇添富北交所创新精选两年定开混合A", "4.95", "65.64", "5.46", "0.2703", 2),
    @(1, "014269", "嘉实北交所精选两年定期混合A", "5.00", "87.43", "4.23", "0.2115", 6),
    @(2, "014294", "南方北交所精选两年定开混合", "4.63", "33.00", "1.44", "0.0667", 5),
    @(3, "014280", "汇添富北交所创新精选两年定开混合C", "0.55", "65.64", "5.46", "0.0300", 2),
    @(4, "014270", "嘉实北交所精选两年定期混合C", "0.64", "87.43", "4.23", "0.0271", 6)
)

for ($i = 0; $i -lt $q1Data.Count; $i++) {
    $row = $i + 2
    $rec = $q1Data[$i]
    $q1.Range("A$row").Value = $rec[0]
    $q1.Range("B$row").Value = "'" + $rec[1]
    $q1.Range("C$row").Value = $rec[2]
    $q1.Range("D$row").Value = "'" + $rec[3]
    $q1.Range("E$row").Value = "'" + $rec[4]
    $q1.Range("F$row").Value = "'" + $rec[5]
    $q1.Range("G$row").Value = "'" + $rec[6]
    $q1.Range("H$row").Value = $rec[7]
}

# ---------------------------------------------------------------------
# 2) New roll-up row in "总计" for 2022-Q1, above the existing rows
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# The blank inserted row inherits stray formatting from the row above;
# clear it so it matches the unstyled data rows elsewhere in the sheet.
$totalSheet.Range("A2:D2").ClearFormats()

$q4Sheet.Range("A2").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 0.61

# The pre-existing rows keep their data, but their running index in column
# A (0,1,2,3) needs to shift to (1,2,3,4) now that a new row-0 is above them.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

# Adding a sheet activates it; restore the original active tab (2021-Q1).
$wb.Worksheets.Item("2021-Q1").Activate()
